$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 80.38461
$ws.Range("I9").Value = 67.5
$ws.Range("J9").Value = 91.42856999999999
$ws.Range("K9").Value = 67.5
$ws.Range("L9").Value = 91.42856999999999
$ws.Range("M9").Value = 101.5
$ws.Range("N9").Value = -429.42857
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H116").Value = 2420.875
$ws.Range("I116").Value = 1866.3889
$ws.Range("J116").Value = 4084.3333
$ws.Range("K116").Value = 1866.3889
$ws.Range("L116").Value = 4084.3333
$ws.Range("M116").Value = 1575.6111
$ws.Range("N116").Value = -10968.3333
$ws.Range("H139").Value = 53000
$ws.Range("J139").Value = 53000
$ws.Range("L139").Value = 53000
$ws.Range("N139").Value = -63280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8622.57
$ws.Range("I32").Value = 6124.852
$ws.Range("J32").Value = 26939.166
$ws.Range("K32").Value = 6124.852
$ws.Range("L32").Value = 26939.166
$ws.Range("M32").Value = -5837.852
$ws.Range("N32").Value = -27513.166
$ws.Range("H132").Value = 714632.8
$ws.Range("I132").Value = 1223320.1
$ws.Range("J132").Value = 6104.0713
$ws.Range("K132").Value = 3669960.3
$ws.Range("L132").Value = 18312.2139
$ws.Range("M132").Value = -3667430.3
$ws.Range("N132").Value = -23372.2139
$ws.Range("H134").Value = 64196
$ws.Range("J134").Value = 64196
$ws.Range("L134").Value = 64196
$ws.Range("N134").Value = -74336
$ws.Range("H139").Value = 30257.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 30257.5
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -40537.5
$ws.Range("H141").Value = 130000
$ws.Range("J141").Value = 130000
$ws.Range("L141").Value = 130000
$ws.Range("N141").Value = -140360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1017.375
$ws.Range("I107").Value = 950.06665
$ws.Range("J107").Value = 1129.5555
$ws.Range("K107").Value = 950.06665
$ws.Range("L107").Value = 1129.5555
$ws.Range("M107").Value = 969.93335
$ws.Range("N107").Value = -4969.5555
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2811.6843
$ws.Range("I16").Value = 2301.4666
$ws.Range("K16").Value = 2301.4666
$ws.Range("M16").Value = -2014.4666
$ws.Range("H62").Value = 7172.1816
$ws.Range("I62").Value = 2527.111
$ws.Range("J62").Value = 28075
$ws.Range("K62").Value = 2527.111
$ws.Range("L62").Value = 28075
$ws.Range("M62").Value = -1903.111
$ws.Range("N62").Value = -29323
$ws.Range("H65").Value = 7172.1816
$ws.Range("I65").Value = 2527.111
$ws.Range("J65").Value = 28075
$ws.Range("K65").Value = 12635.555
$ws.Range("L65").Value = 140375
$ws.Range("M65").Value = -9515.555
$ws.Range("N65").Value = -146615
$ws.Range("H113").Value = 2811.6843
$ws.Range("I113").Value = 2301.4666
$ws.Range("K113").Value = 2301.4666
$ws.Range("M113").Value = -131.4666000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1376.6
$ws.Range("J34").Value = 1629.0834
$ws.Range("L34").Value = 4887.2502
$ws.Range("N34").Value = -5055.2502
$ws.Range("H39").Value = 1043.5
$ws.Range("I39").Value = 937.5
$ws.Range("J39").Value = 1085.9
$ws.Range("K39").Value = 2812.5
$ws.Range("L39").Value = 3257.7
$ws.Range("M39").Value = -2518.5
$ws.Range("N39").Value = -3845.7
$ws.Range("H55").Value = 1396
$ws.Range("I55").Value = 740
$ws.Range("J55").Value = 1833.3334
$ws.Range("K55").Value = 2220
$ws.Range("L55").Value = 5500.0002
$ws.Range("M55").Value = -2043
$ws.Range("N55").Value = -5854.0002
$ws.Range("H86").Value = 1428.5
$ws.Range("I86").Value = 350
$ws.Range("J86").Value = 1644.2
$ws.Range("K86").Value = 1050
$ws.Range("L86").Value = 4932.6
$ws.Range("M86").Value = 136
$ws.Range("N86").Value = -7304.6
$ws.Range("H89").Value = 1428.5
$ws.Range("I89").Value = 350
$ws.Range("J89").Value = 1644.2
$ws.Range("K89").Value = 3150
$ws.Range("L89").Value = 14797.8
$ws.Range("M89").Value = 2778
$ws.Range("N89").Value = -26653.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 13071.429
$ws.Range("J63").Value = 13071.429
$ws.Range("L63").Value = 13071.429
$ws.Range("N63").Value = -14443.429
$ws.Range("H66").Value = 13071.429
$ws.Range("J66").Value = 13071.429
$ws.Range("L66").Value = 39214.287
$ws.Range("N66").Value = -46078.287
$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766
$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652
$ws.Range("H134").Value = 22652.625
$ws.Range("J134").Value = 22652.625
$ws.Range("L134").Value = 67957.875
$ws.Range("N134").Value = -73027.875
$ws.Range("H135").Value = 96391.42999999999
$ws.Range("J135").Value = 96391.42999999999
$ws.Range("L135").Value = 96391.42999999999
$ws.Range("N135").Value = -106531.43
$ws.Range("H141").Value = 19633.334
$ws.Range("J141").Value = 19633.334
$ws.Range("L141").Value = 19633.334
$ws.Range("N141").Value = -29993.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 588838.25
$ws.Range("I46").Value = 440
$ws.Range("J46").Value = 834004.2
$ws.Range("K46").Value = 440
$ws.Range("L46").Value = 834004.2
$ws.Range("M46").Value = -252
$ws.Range("N46").Value = -834380.2
$ws.Range("H95").Value = 18058.545
$ws.Range("J95").Value = 18058.545
$ws.Range("L95").Value = 18058.545
$ws.Range("N95").Value = -23550.545
$ws.Range("H138").Value = 39476.332
$ws.Range("J138").Value = 39476.332
$ws.Range("L138").Value = 39476.332
$ws.Range("N138").Value = -49756.332
$ws.Range("H141").Value = 77211.25
$ws.Range("J141").Value = 77211.25
$ws.Range("L141").Value = 77211.25
$ws.Range("N141").Value = -87571.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H113").Value = 606.63635
$ws.Range("I113").Value = 315.70587
$ws.Range("J113").Value = 1595.8
$ws.Range("K113").Value = 947.11761
$ws.Range("L113").Value = 4787.4
$ws.Range("M113").Value = 1222.88239
$ws.Range("N113").Value = -9127.4
$ws.Range("H140").Value = 34122.582
$ws.Range("J140").Value = 34122.582
$ws.Range("L140").Value = 34122.582
$ws.Range("N140").Value = -44482.582
$ws.Range("H141").Value = 35117.5
$ws.Range("J141").Value = 35117.5
$ws.Range("L141").Value = 35117.5
$ws.Range("N141").Value = -45477.5
